$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole "Tenant Passport ID Number" column (C) and delete it,
# shifting Start date/Finish/... left - mirrors the manual "right-click
# column header > Delete" flow in Excel.
[void]$ws.Columns("C").Select()
[void]$ws.Columns("C").Delete()
